# Adds payment processing background service
# - Adds new "Pay from Credit" related error codes to the Payments sheet
# - Moves the active/selected tab from Auth to Transactions
# - Leaves a C7 selection on the Payments sheet (last place the author worked)

$wb = $excel.ActiveWorkbook

# --- Payments sheet: new error code rows -------------------------------
$wsPayments = $wb.Worksheets.Item("Payments")
$wsPayments.Activate()

$wsPayments.Range("B4").Value = "payment does not exist"
$wsPayments.Range("C4").Value = "ERROR"

$wsPayments.Range("B5").Value = "user does not have permission to pay payment"
$wsPayments.Range("C5").Value = "ERROR"

$wsPayments.Range("B6").Value = "payment is already paid"
$wsPayments.Range("C6").Value = "WARNING"

# Leave the selection where the author last left it on this sheet.
$wsPayments.Range("C7").Select()

# --- Final active sheet: Transactions -----------------------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Activate()
